# Updated symbol list on Tue Dec 20 23:07:16 UTC 2022 with GitHub Actions
# Refreshes prices/hour column and re-ranks a block of exchange-token rows.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '251.38'
$ws.Range("G2").NumberFormat = "@"
$ws.Range("G2").Value = '23'
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '22.83'
$ws.Range("G3").NumberFormat = "@"
$ws.Range("G3").Value = '23'
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '5.424'
$ws.Range("G4").NumberFormat = "@"
$ws.Range("G4").Value = '23'
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.05687'
$ws.Range("G5").NumberFormat = "@"
$ws.Range("G5").Value = '23'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '3.413'
$ws.Range("G6").NumberFormat = "@"
$ws.Range("G6").Value = '23'
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '6.370'
$ws.Range("G7").NumberFormat = "@"
$ws.Range("G7").Value = '23'
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.8134'
$ws.Range("G8").NumberFormat = "@"
$ws.Range("G8").Value = '23'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.9427'
$ws.Range("G9").NumberFormat = "@"
$ws.Range("G9").Value = '23'
$ws.Range("B10").Value = 'WazirX'
$ws.Range("C10").Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.1439'
$ws.Range("E10").Value = '9WazirXWRX'
$ws.Range("G10").NumberFormat = "@"
$ws.Range("G10").Value = '23'
$ws.Range("B11").Value = 'MandalaExchangeToken'
$ws.Range("C11").Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07443'
$ws.Range("E11").Value = '10MandalaExchangeTokenMDX'
$ws.Range("G11").NumberFormat = "@"
$ws.Range("G11").Value = '23'
$ws.Range("B12").Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range("C12").Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.03174'
$ws.Range("E12").Value = '11LiechtensteinCryptoassetsExchangeLCX'
$ws.Range("G12").NumberFormat = "@"
$ws.Range("G12").Value = '23'
$ws.Range("B13").Value = 'BitrueCoin'
$ws.Range("C13").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.03064'
$ws.Range("E13").Value = '12BitrueCoinBTR'
$ws.Range("G13").NumberFormat = "@"
$ws.Range("G13").Value = '23'
$ws.Range("B14").Value = 'BitMartToken'
$ws.Range("C14").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.09368'
$ws.Range("E14").Value = '13BitMartTokenBMX'
$ws.Range("G14").NumberFormat = "@"
$ws.Range("G14").Value = '23'
$ws.Range("B15").Value = 'MCDex'
$ws.Range("C15").Value = 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.718'
$ws.Range("E15").Value = '14MCDexMCB'
$ws.Range("G15").NumberFormat = "@"
$ws.Range("G15").Value = '23'
$ws.Range("B16").Value = 'BitForexToken'
$ws.Range("C16").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.001597'
$ws.Range("E16").Value = '15BitForexTokenBF'
$ws.Range("G16").NumberFormat = "@"
$ws.Range("G16").Value = '23'
$ws.Range("B17").Value = 'CoinExToken'
$ws.Range("C17").Value = 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.04774'
$ws.Range("E17").Value = '16CoinExTokenCET'
$ws.Range("G17").NumberFormat = "@"
$ws.Range("G17").Value = '23'
$ws.Range("B18").Value = 'One'
$ws.Range("C18").Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.0005789'
$ws.Range("E18").Value = '17OneONE'
$ws.Range("G18").NumberFormat = "@"
$ws.Range("G18").Value = '23'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.006378'
$ws.Range("G19").NumberFormat = "@"
$ws.Range("G19").Value = '23'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.005047'
$ws.Range("G20").NumberFormat = "@"
$ws.Range("G20").Value = '23'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.001027'
$ws.Range("G21").NumberFormat = "@"
$ws.Range("G21").Value = '23'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.0001499'
$ws.Range("G22").NumberFormat = "@"
$ws.Range("G22").Value = '23'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '3.719'
$ws.Range("G23").NumberFormat = "@"
$ws.Range("G23").Value = '23'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.181'
$ws.Range("G24").NumberFormat = "@"
$ws.Range("G24").Value = '23'
$ws.Range("G25").NumberFormat = "@"
$ws.Range("G25").Value = '23'
$ws.Range("G26").NumberFormat = "@"
$ws.Range("G26").Value = '23'
$ws.Range("G27").NumberFormat = "@"
$ws.Range("G27").Value = '23'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.0003000'
$ws.Range("G28").NumberFormat = "@"
$ws.Range("G28").Value = '23'
$ws.Range("G29").NumberFormat = "@"
$ws.Range("G29").Value = '23'
$ws.Range("G30").NumberFormat = "@"
$ws.Range("G30").Value = '23'
$ws.Range("G31").NumberFormat = "@"
$ws.Range("G31").Value = '23'
$ws.Range("G32").NumberFormat = "@"
$ws.Range("G32").Value = '23'
$ws.Range("G33").NumberFormat = "@"
$ws.Range("G33").Value = '23'
$ws.Range("G34").NumberFormat = "@"
$ws.Range("G34").Value = '23'
$ws.Range("G35").NumberFormat = "@"
$ws.Range("G35").Value = '23'
$ws.Range("G36").NumberFormat = "@"
$ws.Range("G36").Value = '23'
$ws.Range("G37").NumberFormat = "@"
$ws.Range("G37").Value = '23'
$ws.Range("G38").NumberFormat = "@"
$ws.Range("G38").Value = '23'
$ws.Range("G39").NumberFormat = "@"
$ws.Range("G39").Value = '23'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.04030'
$ws.Range("G40").NumberFormat = "@"
$ws.Range("G40").Value = '23'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.006765'
$ws.Range("G41").NumberFormat = "@"
$ws.Range("G41").Value = '23'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.1073'
$ws.Range("G42").NumberFormat = "@"
$ws.Range("G42").Value = '23'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.002708'
$ws.Range("G43").NumberFormat = "@"
$ws.Range("G43").Value = '23'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.007568'
$ws.Range("G44").NumberFormat = "@"
$ws.Range("G44").Value = '23'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.00005760'
$ws.Range("G45").NumberFormat = "@"
$ws.Range("G45").Value = '23'
$ws.Range("G46").NumberFormat = "@"
$ws.Range("G46").Value = '23'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.5000'
$ws.Range("G47").NumberFormat = "@"
$ws.Range("G47").Value = '23'
$ws.Range("E48").Value = '47BOLOBOLOBestin24h'
$ws.Range("G48").NumberFormat = "@"
$ws.Range("G48").Value = '23'
$ws.Range("G49").NumberFormat = "@"
$ws.Range("G49").Value = '23'
$ws.Range("G50").NumberFormat = "@"
$ws.Range("G50").Value = '23'
$ws.Range("G51").NumberFormat = "@"
$ws.Range("G51").Value = '23'
